$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cell H1 with same formatting as G1 (bold, bordered, centered header style)
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Fill H2:H12 with 0 (plain numeric, no special style - matches other data cells)
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 8).Value = 0
}
